$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 1 above the existing data, shifting everything down.
$ws.Rows.Item(1).Insert()

# Row 1: TotalWorktime
$ws.Range("A1").Value = "TotalWorktime"
$ws.Range("B1").Value = 7298
$ws.Range("C1").Value = 380

# Row 2: TotalProducts
$ws.Range("A2").Value = "TotalProducts"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 235455

# Row 3: TotalGoodProducts
$ws.Range("A3").Value = "TotalGoodProducts"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 222838

# Row 4: TotalScrapProducts
$ws.Range("A4").Value = "TotalScrapProducts"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 12617

# Row 5: MachineSpeed
$ws.Range("A5").Value = "MachineSpeed"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 900

# Row 6: Scrap_Percentage
$ws.Range("A6").Value = "Scrap_Percentage"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 5.36

# Row 7: OEE
$ws.Range("A7").Value = "OEE"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 65.16

# Remove the now-stale rows 8:9 (old row7/row8 remnants shifted down by the insert above).
$ws.Range("A8:A9").EntireRow.Delete()
